# Apply the edits described by the commit "fixes ahead of Solvi":
#  - run_specification gets two new association rows for run_id 3
#    (experiment_id 1 and 2) inserted ahead of the existing rows,
#  - selections/active sheet are updated to reflect where the author
#    was last working (run_specification ends up the active tab),
#  - experiment_specification's scrolled viewport is reset back to A1,
#  - experiment_description's selection moves on one row.

$wb = $excel.ActiveWorkbook

# --- experiment_description: selection moves from C5 to C6 ---
$wsDescription = $wb.Worksheets.Item(1)
$wsDescription.Activate() | Out-Null
$wsDescription.Range("C6").Select() | Out-Null

# --- experiment_specification: scroll position resets to the top-left (A1) ---
$wsSpecification = $wb.Worksheets.Item(2)
$wsSpecification.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$wsSpecification.Range("E46").Select() | Out-Null

# --- run_specification: insert two rows of data for run_id 3 ---
$wsRunSpecification = $wb.Worksheets.Item(4)
$wsRunSpecification.Activate() | Out-Null

$wsRunSpecification.Rows("7:8").Insert() | Out-Null

$wsRunSpecification.Range("A7").Value = 3
$wsRunSpecification.Range("B7").Value = 1

$wsRunSpecification.Range("A8").Value = 3
$wsRunSpecification.Range("B8").Value = 2

# run_specification ends up the active/selected sheet and cell
$wsRunSpecification.Range("C18").Select() | Out-Null
